# Generate Report for handback
# The "defb06e1-2df3-47d9-b3d3-7f40148e592b" file has now been handed back
# (in sync with en-US) for both the zh-cn and de-de locales. Update the
# Overview sheet status and each locale sheet's Status / Latest Handback
# DateTime columns accordingly.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = $statusText
$wsZhCn.Range("G3").Value = "2016-01-26 09:40:07"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = $statusText
$wsDeDe.Range("G3").Value = "2016-01-26 09:40:26"
